$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 need the same style as the existing header cells
# (bold font + thin border + centered/top alignment). Copy the format from
# H1 (an existing header cell) instead of re-building the style by hand so
# the workbook's shared style table stays byte-identical.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-15: new I (I0) and J (IF) numeric columns
$values = @(
    @(2, 8, 8),
    @(3, 7, 7),
    @(4, 8, 8),
    @(5, 8, 8),
    @(6, 7, 7),
    @(7, 7, 7),
    @(8, 2, 3),
    @(9, 5, 5),
    @(10, 9, 9),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 8, 9),
    @(14, 4, 4),
    @(15, 3, 3)
)

foreach ($row in $values) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
